$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (row 239 through row 244), columns A (date serial), B, C, D
$data = @(
    @(239, 44313, 2, 13, 81.40262993112086),
    @(240, 44314, 2, 15, 93.9261114589856),
    @(241, 44315, 5, 19, 118.9730745147151),
    @(242, 44316, 1, 16, 100.187852222918),
    @(243, 44317, 5, 20, 125.2348152786475),
    @(244, 44318, 2, 19, 118.9730745147151)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r - 1, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
